$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Column A width: 56.7109375 -> ~63.140625 chars (closest achievable: 63.1667)
$ws.Columns.Item(1).ColumnWidth = 62.3

# 2. Row 2: 'Statistika igre' description moves to column A; X marks in B/C/D
$ws.Range("A2").Value = "Statistika igre (Ukupno prodaja i ukupna zarada po državama i po vremenu)"
$ws.Range("B2").Value = "X"
$ws.Range("C2").Value = "X"
$ws.Range("D2").Value = "X"

# 3. Row 3: 'Statistika developera' description (text updated); X marks in B/D/E
$ws.Range("A3").Value = "Statistika developera (Ocena najbolje igre, prosečna cena igrica, prosečna ocena igrica, ukupan broj prodatih igara najboljih developera)"
$ws.Range("B3").Value = "X"
$ws.Range("D3").Value = "X"
$ws.Range("E3").Value = "X"

# 4. Row 4: 'Statistika po državi' description (text updated); X marks in B/C/D/F
$ws.Range("A4").Value = "Statistika po državi (Prosečna količina novca u wallet-u po državi, ukupno prodato igara po državi, suma novca potrošena po državi, prosečna ukupna cena kupovine, prodaja po žanrovima po državi, prodaja po državi tokom vremena sa informacijama o wallet-u)"
$ws.Range("B4").Value = "X"
$ws.Range("C4").Value = "X"
$ws.Range("D4").Value = "X"
$ws.Range("F4").Value = "X"

# 5. Row 5 is a brand-new data row (previously blank) - clone formatting from row 3
$ws.Range("A3:Y3").Copy()
$ws.Range("A5:Y5").PasteSpecial(-4122)
$ws.Range("A5:Y5").RowHeight = 45.75
$ws.Range("A5").Value = "Statistika žanra (prodaja po žanrovima po državi, najprodavaniji žanrovi, statistika 3 najprodavanija žanra po državi)"
$ws.Range("B5").Value = "X"
$ws.Range("C5").Value = "X"
$ws.Range("D5").Value = "X"

# 6. Grow the Table_1 ListObject to cover the new row
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:Y5"))

# 7. Split the single color-scale conditional format into two rules:
#    - existing A1:Y1/A3:Y4 rule drops to priority 2
#    - new A5:Y5 rule takes priority 1, same green->white color scale
$newCf = $ws.Range("A5:Y5").FormatConditions.AddColorScale(2)
$newCf.ColorScaleCriteria(1).FormatColor.Color = 9091927
$newCf.ColorScaleCriteria(2).FormatColor.Color = 16777215
$newCf.Priority = 1
$existingCf = $ws.Range("A1:Y1").FormatConditions.Item(1)
$existingCf.Priority = 2

# 8. Restore the active selection at E6
$ws.Range("E6").Select()
